$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{ Row=2; Price="67.707.43"; Vol="  -1.61%  " },
    @{ Row=3; Price="3.268.79"; Vol="  -0.46%  " },
    @{ Row=4; Price=$null; Vol="  -0.06%  " },
    @{ Row=5; Price="580.58"; Vol="  -0.63%  " },
    @{ Row=6; Price="184.64"; Vol="  +1.29%  " },
    @{ Row=7; Price=$null; Vol="  -0.02%  " },
    @{ Row=8; Price=$null; Vol="  +0.77%  " },
    @{ Row=9; Price=$null; Vol="  -3.16%  " },
    @{ Row=10; Price="6.54"; Vol="  -2.03%  " },
    @{ Row=11; Price="0.408"; Vol="  -4.58%  " },
    @{ Row=12; Price="3.834.88"; Vol="  -0.43%  " },
    @{ Row=13; Price=$null; Vol="  +0.65%  " },
    @{ Row=14; Price="27.37"; Vol="  -5.27%  " },
    @{ Row=15; Price="67.736.59"; Vol="  -1.61%  " },
    @{ Row=16; Price=$null; Vol="  -2.25%  " },
    @{ Row=17; Price="3.252.30"; Vol="  -0.85%  " },
    @{ Row=18; Price="5.70"; Vol="  -2.62%  " },
    @{ Row=19; Price="13.42"; Vol="  -1.42%  " },
    @{ Row=20; Price="401.19"; Vol="  +1.45%  " },
    @{ Row=21; Price=$null; Vol="  -2.18%  " },
    @{ Row=22; Price=$null; Vol="  +0.10%  " },
    @{ Row=23; Price=$null; Vol="  -1.59%  " },
    @{ Row=24; Price="0.507"; Vol="  -1.93%  " },
    @{ Row=25; Price="0.0000117"; Vol="  -2.42%  " },
    @{ Row=26; Price="0.186"; Vol="  -1.77%  " },
    @{ Row=27; Price="9.49"; Vol="  -2.43%  " },
    @{ Row=28; Price=$null; Vol="  +0.37%  " },
    @{ Row=29; Price="1.95"; Vol="  -1.88%  " },
    @{ Row=30; Price="22.62"; Vol="  -1.80%  " },
    @{ Row=31; Price=$null; Vol="  -4.33%  " },
    @{ Row=32; Price="6.91"; Vol="  -3.51%  " },
    @{ Row=33; Price=$null; Vol="  +0.07%  " },
    @{ Row=34; Price="1.24"; Vol="  -4.30%  " },
    @{ Row=35; Price="164.24"; Vol="  -0.22%  " },
    @{ Row=36; Price=$null; Vol="  -4.02%  " },
    @{ Row=37; Price="1.88"; Vol="  -2.02%  " },
    @{ Row=38; Price="26.92"; Vol="  +1.93%  " },
    @{ Row=39; Price=$null; Vol="  -3.51%  " },
    @{ Row=40; Price=$null; Vol="  -2.13%  " },
    @{ Row=41; Price="6.36"; Vol="  -3.55%  " },
    @{ Row=42; Price="2.681.00"; Vol="  +2.04%  " },
    @{ Row=45; Price="0.0677"; Vol="  -1.77%  " },
    @{ Row=48; Price=$null; Vol="  -3.42%  " },
    @{ Row=49; Price=$null; Vol="  -0.77%  " },
    @{ Row=50; Price=$null; Vol="  -1.98%  " },
    @{ Row=51; Price=$null; Vol="  -2.45%  " },
)

foreach ($item in $data) {
    if ($item.Price -ne $null) {
        $ws.Range("D$($item.Row)").NumberFormat = "@"
        $ws.Range("D$($item.Row)").Value = $item.Price
    }
    $ws.Range("E$($item.Row)").NumberFormat = "@"
    $ws.Range("E$($item.Row)").Value = $item.Vol
}

# Rows 43/44: OKB and dogwifhat swap places (with slightly updated values)
$ws.Range("B43").Value = "OKB"
$ws.Range("C43").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "40.71"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -1.77%  "

$ws.Range("B44").Value = "dogwifhat"
$ws.Range("C44").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.44"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -4.34%  "

# Rows 46/47: Bittensor and InjectiveProtocol swap places (with slightly updated values)
$ws.Range("B46").Value = "Bittensor"
$ws.Range("C46").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "334.72"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -3.27%  "

$ws.Range("B47").Value = "InjectiveProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "24.53"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.73%  "
